$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row for Netherlands, alphabetically between Luxembourg
#    (row 20) and New Zealand (old row 21 -> shifts down to row 22).
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).Insert()

# ---------------------------------------------------------------------------
# 2) Insert a new row for Taiwan, alphabetically between Switzerland
#    (row 28 after the first insert) and UK (old row 28 -> shifts to row 29,
#    then to row 30 after this insert).
# ---------------------------------------------------------------------------
$ws.Rows.Item(29).Insert()

# ---------------------------------------------------------------------------
# Reset styling across the whole table back to the workbook default ("Normal"
# cell style) and reapply the percentage format to the %Excess column, same
# as the re-saved workbook does.
# ---------------------------------------------------------------------------
$ws.Range("A1:G31").Style = "Normal"
$ws.Range("F2:F31").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# Netherlands data (row 21)
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Netherlands"
$ws.Range("B21").Value = "17.4M"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "325,475"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "339,650"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "14,175"

$ws.Range("F21").Value = 0.044
$ws.Range("G21").Value = "(+1.9%, +7.9%)"

# ---------------------------------------------------------------------------
# Taiwan data (row 29)
# ---------------------------------------------------------------------------
$ws.Range("A29").Value = "Taiwan"
$ws.Range("B29").Value = "23.6M"

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "383,471"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "357,239"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "-26,232"

$ws.Range("F29").Value = -0.068
$ws.Range("G29").Value = "(-10.4%, -3.2%)"

# ---------------------------------------------------------------------------
# Match the saved selection in the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("G29").Select()
